$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.759.39"
$ws.Range("E2").Value = "  +5.82%  "
$ws.Range("D3").Value = "2.671.45"
$ws.Range("E3").Value = "  +7.01%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.05"
$ws.Range("E5").Value = "  +3.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.21"
$ws.Range("E6").Value = "  +5.10%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +3.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.177"
$ws.Range("E9").Value = "  +17.67%  "
$ws.Range("D10").Value = "2.676.62"
$ws.Range("E10").Value = "  +7.21%  "
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.353"
$ws.Range("E12").Value = "  +6.04%  "
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("E14").Value = "  +13.87%  "
$ws.Range("D15").Value = "3.121.41"
$ws.Range("E15").Value = "  +5.74%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "72.713.01"
$ws.Range("E16").Value = "  +6.05%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.94"
$ws.Range("E17").Value = "  +5.37%  "
$ws.Range("D18").Value = "2.674.16"
$ws.Range("E18").Value = "  +7.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "386.74"
$ws.Range("E19").Value = "  +7.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.67"
$ws.Range("E20").Value = "  +7.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.96"
$ws.Range("E21").Value = "  +6.98%  "
$ws.Range("E22").Value = "  +5.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.05"
$ws.Range("E23").Value = "  +25.97%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.68"
$ws.Range("E24").Value = "  +5.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.48"
$ws.Range("E25").Value = "  +8.92%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.04"
$ws.Range("E27").Value = "  +13.99%  "
$ws.Range("D28").Value = "2.815.44"
$ws.Range("E28").Value = "  +6.94%  "
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "0.0₃0983"
$ws.Range("E30").Value = "  +13.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "550.77"
$ws.Range("E31").Value = "  +9.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  +6.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  +12.80%  "
$ws.Range("E34").Value = "  +5.43%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.07"
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.43"
$ws.Range("E37").Value = "  +4.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.115"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("E39").Value = "  +10.68%  "
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.86"
$ws.Range("E41").Value = "  +10.96%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.14"
$ws.Range("E42").Value = "  +9.54%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.65"
$ws.Range("E43").Value = "  +16.50%  "
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.337"
$ws.Range("E45").Value = "  +7.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.84"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "153.12"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.71"
$ws.Range("E48").Value = "  +5.64%  "
$ws.Range("D49").Value = "0.0₆0273"
$ws.Range("E49").Value = "  +13.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.546"
$ws.Range("E50").Value = "  +7.72%  "
$ws.Range("E51").Value = "  +11.38%  "
